$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 16.48560822670832
$ws.Range("C2").Value = 13.14420653597452
$ws.Range("D2").Value = 4.269846445367235
$ws.Range("F2").Value = 20.41721901048637
$ws.Range("G2").Value = 22.10947758239378
$ws.Range("H2").Value = 12.61642110207039
$ws.Range("L2").Value = 10.64346660949809
$ws.Range("O2").Value = 18.31204254687029
$ws.Range("B3").Value = 15.74963793196334
$ws.Range("C3").Value = 13.02269654296349
$ws.Range("D3").Value = 4.19332347410278
$ws.Range("F3").Value = 20.47597134890049
$ws.Range("G3").Value = 22.19999390981196
$ws.Range("H3").Value = 12.67831197238196
$ws.Range("L3").Value = 10.61508727954655
$ws.Range("O3").Value = 18.41327403191698
$ws.Range("B4").Value = 15.27978942141904
$ws.Range("C4").Value = 12.94900655043707
$ws.Range("D4").Value = 4.145233317719673
$ws.Range("F4").Value = 20.51990974962818
$ws.Range("G4").Value = 22.2675518785262
$ws.Range("H4").Value = 12.71908081149453
$ws.Range("L4").Value = 10.59981703790943
$ws.Range("O4").Value = 18.4812157299647
$ws.Range("B5").Value = 15.08403621133034
$ws.Range("C5").Value = 12.91923511180788
$ws.Range("D5").Value = 4.125374280529976
$ws.Range("F5").Value = 20.53978218561496
$ws.Range("G5").Value = 22.2980681024273
$ws.Range("H5").Value = 12.7363896453358
$ws.Range("L5").Value = 10.59414064662122
$ws.Range("O5").Value = 18.5103510506312
$ws.Range("B6").Value = 15.05128005296383
$ws.Range("C6").Value = 12.91430794362724
$ws.Range("D6").Value = 4.122061360668447
$ws.Range("F6").Value = 20.54320048032119
$ws.Range("G6").Value = 22.3033148126474
$ws.Range("H6").Value = 12.73930573141857
$ws.Range("L6").Value = 10.59323120640068
$ws.Range("O6").Value = 18.51527625404389
$ws.Range("B7").Value = 15.27716645095815
$ws.Range("C7").Value = 12.94860396296135
$ws.Range("D7").Value = 4.144966530806118
$ws.Range("F7").Value = 20.52016980560385
$ws.Range("G7").Value = 22.26795137845429
$ws.Range("H7").Value = 12.71931143060232
$ws.Range("L7").Value = 10.59973826635853
$ws.Range("O7").Value = 18.48160280229223
$ws.Range("B8").Value = 16.23570107714414
$ws.Range("C8").Value = 13.10213578577606
$ws.Range("D8").Value = 4.243699943344268
$ws.Range("F8").Value = 20.43583866973199
$ws.Range("G8").Value = 22.13818551777102
$ws.Range("H8").Value = 12.63718621640669
$ws.Range("L8").Value = 10.63323705862205
$ws.Range("O8").Value = 18.34574274090117
$ws.Range("B9").Value = 17.96437728297153
$ws.Range("C9").Value = 13.40928987386695
$ws.Range("D9").Value = 4.427887028261143
$ws.Range("F9").Value = 20.33328076712987
$ws.Range("G9").Value = 21.97987289322772
$ws.Range("H9").Value = 12.49813149256326
$ws.Range("L9").Value = 10.71580401432846
$ws.Range("O9").Value = 18.12549918310982
$ws.Range("B10").Value = 19.13295645111236
$ws.Range("C10").Value = 13.63704245621771
$ws.Range("D10").Value = 4.556624795123002
$ws.Range("F10").Value = 20.29672312591429
$ws.Range("G10").Value = 21.92352105270185
$ws.Range("H10").Value = 12.40942067038008
$ws.Range("L10").Value = 10.78642520074493
$ws.Range("O10").Value = 17.99221011865682
$ws.Range("B11").Value = 19.64103167021725
$ws.Range("C11").Value = 13.7407623821138
$ws.Range("D11").Value = 4.613599187904573
$ws.Range("F11").Value = 20.2885976310612
$ws.Range("G11").Value = 21.91113883543983
$ws.Range("H11").Value = 12.37199521777546
$ws.Range("L11").Value = 10.82063482835771
$ws.Range("O11").Value = 17.93784591910553
$ws.Range("B12").Value = 19.82994552696034
$ws.Range("C12").Value = 13.78002625349406
$ws.Range("D12").Value = 4.634932878021933
$ws.Range("F12").Value = 20.28674845102175
$ws.Range("G12").Value = 21.90837071397378
$ws.Range("H12").Value = 12.35824539936159
$ws.Range("L12").Value = 10.83388121672248
$ws.Range("O12").Value = 17.91816761222406
$ws.Range("B13").Value = 19.78941590028172
$ws.Range("C13").Value = 13.7715710951149
$ws.Range("D13").Value = 4.630349211762272
$ws.Range("F13").Value = 20.28709203744933
$ws.Range("G13").Value = 21.90888124424453
$ws.Range("H13").Value = 12.36118786479451
$ws.Range("L13").Value = 10.83101551270869
$ws.Range("O13").Value = 17.9223651885765
$ws.Range("B14").Value = 19.6566440563719
$ws.Range("C14").Value = 13.74399306206677
$ws.Range("D14").Value = 4.615359229255097
$ws.Range("F14").Value = 20.28842086883958
$ws.Range("G14").Value = 21.91087252993233
$ws.Range("H14").Value = 12.37085554012937
$ws.Range("L14").Value = 10.82171881128051
$ws.Range("O14").Value = 17.93620873272608
$ws.Range("B15").Value = 19.57486113223653
$ws.Range("C15").Value = 13.72709822318894
$ws.Range("D15").Value = 4.606145645578152
$ws.Range("F15").Value = 20.28939482786783
$ws.Range("G15").Value = 21.91234278525461
$ws.Range("H15").Value = 12.37683231159149
$ws.Range("L15").Value = 10.81606209740294
$ws.Range("O15").Value = 17.94480677850904
$ws.Range("B16").Value = 19.09927024572125
$ws.Range("C16").Value = 13.63026392191507
$ws.Range("D16").Value = 4.552868298571943
$ws.Range("F16").Value = 20.29742572352128
$ws.Range("G16").Value = 21.92459828896759
$ws.Range("H16").Value = 12.41192555711747
$ws.Range("L16").Value = 10.78423080332345
$ws.Range("O16").Value = 17.99588968056576
$ws.Range("B17").Value = 18.80140954565891
$ws.Range("C17").Value = 13.57086767886251
$ws.Range("D17").Value = 4.519768033859395
$ws.Range("F17").Value = 20.30453451859366
$ws.Range("G17").Value = 21.93552257806825
$ws.Range("H17").Value = 12.4342052447783
$ws.Range("L17").Value = 10.76523169146864
$ws.Range("O17").Value = 18.02883796986348
$ws.Range("B18").Value = 18.62787926468965
$ws.Range("C18").Value = 13.53671698998981
$ws.Range("D18").Value = 4.500580779670303
$ws.Range("F18").Value = 20.30942361309434
$ws.Range("G18").Value = 21.94305261161799
$ws.Range("H18").Value = 12.44729560574974
$ws.Range("L18").Value = 10.75450052126731
$ws.Range("O18").Value = 18.04837869172579
$ws.Range("B19").Value = 18.56874892594861
$ws.Range("C19").Value = 13.52515718392328
$ws.Range("D19").Value = 4.494059140946462
$ws.Range("F19").Value = 20.31121626181822
$ws.Range("G19").Value = 21.9458157615887
$ws.Range("H19").Value = 12.45177509050719
$ws.Range("L19").Value = 10.75090112745267
$ws.Range("O19").Value = 18.0550959277483
$ws.Range("B20").Value = 18.83334670545946
$ws.Range("C20").Value = 13.57718943232286
$ws.Range("D20").Value = 4.52330711476987
$ws.Range("F20").Value = 20.30369490787388
$ws.Range("G20").Value = 21.93423052453667
$ws.Range("H20").Value = 12.43180499470789
$ws.Range("L20").Value = 10.76723388458837
$ws.Range("O20").Value = 18.02526948401981
$ws.Range("B21").Value = 19.69573762301156
$ws.Range("C21").Value = 13.75209397053217
$ws.Range("D21").Value = 4.61976879643507
$ws.Range("F21").Value = 20.28799720591589
$ws.Range("G21").Value = 21.91023540703102
$ws.Range("H21").Value = 12.36800443696623
$ws.Range("L21").Value = 10.82444161798283
$ws.Range("O21").Value = 17.93211784634183
$ws.Range("B22").Value = 20.23902159599027
$ws.Range("C22").Value = 13.86631751835618
$ws.Range("D22").Value = 4.68140038260847
$ws.Range("F22").Value = 20.28489550662889
$ws.Range("G22").Value = 21.90575370343856
$ws.Range("H22").Value = 12.32876962472159
$ws.Range("L22").Value = 10.86352798844736
$ws.Range("O22").Value = 17.87653477960297
$ws.Range("B23").Value = 19.95095091557137
$ws.Range("C23").Value = 13.80537167256642
$ws.Range("D23").Value = 4.648639627525849
$ws.Range("F23").Value = 20.28589477568115
$ws.Range("G23").Value = 21.90711662625742
$ws.Range("H23").Value = 12.34948428238753
$ws.Range("L23").Value = 10.84251415856439
$ws.Range("O23").Value = 17.90571366403359
$ws.Range("B24").Value = 18.8189150286555
$ws.Range("C24").Value = 13.57433137468799
$ws.Range("D24").Value = 4.521707585797827
$ws.Range("F24").Value = 20.30407199760102
$ws.Range("G24").Value = 21.93481077096476
$ws.Range("H24").Value = 12.43288927106605
$ws.Range("L24").Value = 10.76632809546731
$ws.Range("O24").Value = 18.02688093325222
$ws.Range("B25").Value = 17.51396238774664
$ws.Range("C25").Value = 13.32571282970673
$ws.Range("D25").Value = 4.379158222087086
$ws.Range("F25").Value = 20.35424471368938
$ws.Range("G25").Value = 22.01225844959903
$ws.Range("H25").Value = 12.53339042960087
$ws.Range("L25").Value = 10.69169270539574
$ws.Range("O25").Value = 18.18009856426754
